$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.741.03"
$ws.Range("E2").Value = "  -1.75%  "

$ws.Range("D3").Value = "1.762.72"
$ws.Range("E3").Value = "  -1.83%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.97"
$ws.Range("E5").Value = "  -1.85%  "

$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4439"
$ws.Range("E7").Value = "  -2.15%  "

$ws.Range("E8").Value = "  +0.77%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.51"
$ws.Range("E9").Value = "  +1.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07799"
$ws.Range("E10").Value = "  +3.25%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.130"
$ws.Range("E11").Value = "  -1.01%  "

$ws.Range("E12").Value = "  +0.10%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.80"
$ws.Range("E13").Value = "  -2.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.207"
$ws.Range("E14").Value = "  -1.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.389"
$ws.Range("E15").Value = "  -1.02%  "

$ws.Range("D16").Value = "1.759.51"
$ws.Range("E16").Value = "  -1.96%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.38"
$ws.Range("E17").Value = "  +13.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001086"
$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06258"
$ws.Range("E19").Value = "  -7.15%  "

$ws.Range("E20").Value = "  +0.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.40"
$ws.Range("E21").Value = "  -0.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.196"
$ws.Range("E22").Value = "  -2.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5332"
$ws.Range("E23").Value = "  -2.45%  "

$ws.Range("D24").Value = "27.778.95"
$ws.Range("E24").Value = "  -1.52%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.67"
$ws.Range("E25").Value = "  -1.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.328"
$ws.Range("E26").Value = "  -3.48%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.94"
$ws.Range("E27").Value = "  +2.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.93"
$ws.Range("E28").Value = "  +1.23%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.374"
$ws.Range("E29").Value = "  +0.55%  "

$ws.Range("D30").Value = "1.957.66"
$ws.Range("E30").Value = "  -1.77%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "129.58"
$ws.Range("E31").Value = "  -2.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.220"
$ws.Range("E32").Value = "  -1.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.809"
$ws.Range("E33").Value = "  +0.27%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09310"
$ws.Range("E34").Value = "  -1.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.651"
$ws.Range("E35").Value = "  -9.63%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.80"
$ws.Range("E36").Value = "  +5.91%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2191"
$ws.Range("E37").Value = "  -8.72%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02334"
$ws.Range("E38").Value = "  -0.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06171"
$ws.Range("E39").Value = "  -2.16%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6518"
$ws.Range("E40").Value = "  -0.53%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.109"
$ws.Range("E41").Value = "  -1.86%  "

$ws.Range("E42").Value = "  -0.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.043"
$ws.Range("E43").Value = "  -3.48%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.419"
$ws.Range("E44").Value = "  -4.08%  "

$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.94"
$ws.Range("E46").Value = "  -1.57%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6037"
$ws.Range("E47").Value = "  -0.67%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.771"
$ws.Range("E48").Value = "  -1.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "126.48"
$ws.Range("E49").Value = "  -2.60%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.003"
$ws.Range("E50").Value = "  -1.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.147"
$ws.Range("E51").Value = "  -0.82%  "
